# Edit: prepend one new week of "Pimiento" (Agricola del Norte S.A. de Arica) price
# records to the daily logic sheet. This shifts the existing data block (rows 888-952)
# down by 6 rows (to 894-958) and fills the freed rows 888-893 with the new week's data
# (fecha serial 45021).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new blank rows right before the existing data block that starts at row 888.
# This pushes rows 888:952 down to 894:958 and keeps their formatting/styles intact,
# picking up the row-888 formatting (including the date style on column D) for the
# newly inserted rows.
$ws.Range("A888:R893").EntireRow.Insert()

# New row 888 - Zafiro rojo / Primera
$ws.Range("A888").Value = 1
$ws.Range("B888").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C888").Value = "Arica y Parinacota"
$ws.Range("D888").Value = 45021
$ws.Range("E888").Value = 15
$ws.Range("F888").Value = 100112002
$ws.Range("G888").Value = "Pimiento"
$ws.Range("H888").Value = "Zafiro rojo"
$ws.Range("I888").Value = "Primera"
$ws.Range("J888").Value = 120
$ws.Range("K888").Value = 14000
$ws.Range("L888").Value = 15000
$ws.Range("M888").Value = 14500
$ws.Range("N888").Value = "`$/caja 15 kilos"
$ws.Range("O888").Value = "Región de Arica y Parinacota"
$ws.Range("P888").Value = 967
$ws.Range("Q888").Value = 15
$ws.Range("R888").Value = "Hortaliza"

# New row 889 - Zafiro rojo / Segunda
$ws.Range("A889").Value = 1
$ws.Range("B889").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C889").Value = "Arica y Parinacota"
$ws.Range("D889").Value = 45021
$ws.Range("E889").Value = 15
$ws.Range("F889").Value = 100112002
$ws.Range("G889").Value = "Pimiento"
$ws.Range("H889").Value = "Zafiro rojo"
$ws.Range("I889").Value = "Segunda"
$ws.Range("J889").Value = 140
$ws.Range("K889").Value = 12000
$ws.Range("L889").Value = 13000
$ws.Range("M889").Value = 12500
$ws.Range("N889").Value = "`$/caja 15 kilos"
$ws.Range("O889").Value = "Región de Arica y Parinacota"
$ws.Range("P889").Value = 833
$ws.Range("Q889").Value = 15
$ws.Range("R889").Value = "Hortaliza"

# New row 890 - Zafiro rojo / Tercera
$ws.Range("A890").Value = 1
$ws.Range("B890").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C890").Value = "Arica y Parinacota"
$ws.Range("D890").Value = 45021
$ws.Range("E890").Value = 15
$ws.Range("F890").Value = 100112002
$ws.Range("G890").Value = "Pimiento"
$ws.Range("H890").Value = "Zafiro rojo"
$ws.Range("I890").Value = "Tercera"
$ws.Range("J890").Value = 140
$ws.Range("K890").Value = 10000
$ws.Range("L890").Value = 11000
$ws.Range("M890").Value = 10500
$ws.Range("N890").Value = "`$/caja 15 kilos"
$ws.Range("O890").Value = "Región de Arica y Parinacota"
$ws.Range("P890").Value = 700
$ws.Range("Q890").Value = 15
$ws.Range("R890").Value = "Hortaliza"

# New row 891 - Zafiro verde / Primera
$ws.Range("A891").Value = 1
$ws.Range("B891").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C891").Value = "Arica y Parinacota"
$ws.Range("D891").Value = 45021
$ws.Range("E891").Value = 15
$ws.Range("F891").Value = 100112002
$ws.Range("G891").Value = "Pimiento"
$ws.Range("H891").Value = "Zafiro verde"
$ws.Range("I891").Value = "Primera"
$ws.Range("J891").Value = 100
$ws.Range("K891").Value = 9000
$ws.Range("L891").Value = 10000
$ws.Range("M891").Value = 9500
$ws.Range("N891").Value = "`$/caja 15 kilos"
$ws.Range("O891").Value = "Región de Arica y Parinacota"
$ws.Range("P891").Value = 633
$ws.Range("Q891").Value = 15
$ws.Range("R891").Value = "Hortaliza"

# New row 892 - Zafiro verde / Segunda
$ws.Range("A892").Value = 1
$ws.Range("B892").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C892").Value = "Arica y Parinacota"
$ws.Range("D892").Value = 45021
$ws.Range("E892").Value = 15
$ws.Range("F892").Value = 100112002
$ws.Range("G892").Value = "Pimiento"
$ws.Range("H892").Value = "Zafiro verde"
$ws.Range("I892").Value = "Segunda"
$ws.Range("J892").Value = 120
$ws.Range("K892").Value = 7000
$ws.Range("L892").Value = 8000
$ws.Range("M892").Value = 7500
$ws.Range("N892").Value = "`$/caja 15 kilos"
$ws.Range("O892").Value = "Región de Arica y Parinacota"
$ws.Range("P892").Value = 500
$ws.Range("Q892").Value = 15
$ws.Range("R892").Value = "Hortaliza"

# New row 893 - Zafiro verde / Tercera
$ws.Range("A893").Value = 1
$ws.Range("B893").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C893").Value = "Arica y Parinacota"
$ws.Range("D893").Value = 45021
$ws.Range("E893").Value = 15
$ws.Range("F893").Value = 100112002
$ws.Range("G893").Value = "Pimiento"
$ws.Range("H893").Value = "Zafiro verde"
$ws.Range("I893").Value = "Tercera"
$ws.Range("J893").Value = 120
$ws.Range("K893").Value = 5000
$ws.Range("L893").Value = 6000
$ws.Range("M893").Value = 5500
$ws.Range("N893").Value = "`$/caja 15 kilos"
$ws.Range("O893").Value = "Región de Arica y Parinacota"
$ws.Range("P893").Value = 367
$ws.Range("Q893").Value = 15
$ws.Range("R893").Value = "Hortaliza"
